$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "29.870.21"
Set-TextValue $ws.Range("E2") "  +0.10%  "
Set-TextValue $ws.Range("D3") "1.887.57"
Set-TextValue $ws.Range("E3") "  -0.20%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D6") "242.61"
Set-TextValue $ws.Range("E6") "  -0.91%  "
Set-TextValue $ws.Range("E7") "  -0.05%  "
Set-TextValue $ws.Range("D8") "0.3123"
Set-TextValue $ws.Range("E8") "  -0.72%  "
Set-TextValue $ws.Range("E9") "  +0.84%  "
Set-TextValue $ws.Range("D10") "0.07167"
Set-TextValue $ws.Range("E10") "  -4.62%  "
Set-TextValue $ws.Range("D11") "0.08566"
Set-TextValue $ws.Range("E11") "  +5.43%  "
Set-TextValue $ws.Range("D12") "0.7633"
Set-TextValue $ws.Range("E12") "  -0.62%  "
Set-TextValue $ws.Range("D13") "1.913.95"
Set-TextValue $ws.Range("E13") "  +0.71%  "
Set-TextValue $ws.Range("D14") "5.363"
Set-TextValue $ws.Range("E14") "  -1.74%  "
Set-TextValue $ws.Range("D15") "93.60"
Set-TextValue $ws.Range("E15") "  +1.40%  "
Set-TextValue $ws.Range("D16") "6.153"
Set-TextValue $ws.Range("E16") "  -0.88%  "
Set-TextValue $ws.Range("D17") "29.861.18"
Set-TextValue $ws.Range("E17") "  +0.00%  "
Set-TextValue $ws.Range("E18") "  -1.58%  "
Set-TextValue $ws.Range("D19") "244.49"
Set-TextValue $ws.Range("E19") "  +0.02%  "
Set-TextValue $ws.Range("D20") "0.000007806"
Set-TextValue $ws.Range("E20") "  -1.19%  "
Set-TextValue $ws.Range("D21") "2.144.12"
Set-TextValue $ws.Range("E21") "  -0.41%  "
Set-TextValue $ws.Range("E22") "  -0.21%  "
Set-TextValue $ws.Range("D23") "8.016"
Set-TextValue $ws.Range("E23") "  -1.02%  "
Set-TextValue $ws.Range("E24") "  +0.03%  "
Set-TextValue $ws.Range("D25") "0.1638"
Set-TextValue $ws.Range("E25") "  +4.05%  "
Set-TextValue $ws.Range("D26") "9.386"
Set-TextValue $ws.Range("E26") "  -0.60%  "
Set-TextValue $ws.Range("D27") "163.05"
Set-TextValue $ws.Range("E27") "  +0.11%  "
Set-TextValue $ws.Range("D28") "18.73"
Set-TextValue $ws.Range("D29") "2.031"
Set-TextValue $ws.Range("E29") "  -0.62%  "
Set-TextValue $ws.Range("D30") "1.468"
Set-TextValue $ws.Range("E30") "  +2.37%  "
Set-TextValue $ws.Range("D31") "1.539"
Set-TextValue $ws.Range("E31") "  -0.83%  "
Set-TextValue $ws.Range("D32") "4.507"
Set-TextValue $ws.Range("E32") "  +0.12%  "
Set-TextValue $ws.Range("D33") "4.096"
Set-TextValue $ws.Range("E33") "  -0.08%  "
Set-TextValue $ws.Range("D34") "0.05447"
Set-TextValue $ws.Range("E34") "  -1.34%  "
Set-TextValue $ws.Range("E35") "  -1.07%  "
Set-TextValue $ws.Range("D36") "0.7425"
Set-TextValue $ws.Range("E36") "  -2.13%  "
Set-TextValue $ws.Range("D37") "1.001"
Set-TextValue $ws.Range("E37") "  -0.16%  "
Set-TextValue $ws.Range("D38") "2.694"
Set-TextValue $ws.Range("E38") "  +1.94%  "
Set-TextValue $ws.Range("D39") "0.01952"
Set-TextValue $ws.Range("E39") "  +1.47%  "
Set-TextValue $ws.Range("D40") "2.784"
Set-TextValue $ws.Range("E40") "  -0.12%  "
Set-TextValue $ws.Range("D41") "0.4468"
Set-TextValue $ws.Range("E41") "  +0.24%  "
Set-TextValue $ws.Range("D42") "1.105.36"
Set-TextValue $ws.Range("E42") "  -5.31%  "
Set-TextValue $ws.Range("D43") "73.21"
Set-TextValue $ws.Range("E43") "  -1.04%  "
Set-TextValue $ws.Range("D44") "6.070"
Set-TextValue $ws.Range("E44") "  +1.73%  "
Set-TextValue $ws.Range("D45") "0.8518"
Set-TextValue $ws.Range("E45") "  +0.45%  "
Set-TextValue $ws.Range("D46") "1.000"
Set-TextValue $ws.Range("E46") "  -0.06%  "
Set-TextValue $ws.Range("D47") "102.64"
Set-TextValue $ws.Range("E47") "  +0.38%  "
Set-TextValue $ws.Range("D48") "7.661"
Set-TextValue $ws.Range("E48") "  +1.23%  "
Set-TextValue $ws.Range("E49") "  -2.33%  "
Set-TextValue $ws.Range("D50") "3.010"
Set-TextValue $ws.Range("E50") "  -3.00%  "
Set-TextValue $ws.Range("D51") "2.044.69"
Set-TextValue $ws.Range("E51") "  +0.34%  "
